$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.917.93'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '1.632.29'
$ws.Range('E3').Value = '  -2.64%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''209.59'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').Value = '''0.5200'
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''0.2565'
$ws.Range('E8').Value = '  -3.42%  '
$ws.Range('D9').Value = '''0.06231'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('E10').Value = '  -5.30%  '
$ws.Range('D11').Value = '''0.07571'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.625.60'
$ws.Range('E12').Value = '  -2.93%  '
$ws.Range('D13').Value = '''4.350'
$ws.Range('D14').Value = '1.860.60'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '''0.5417'
$ws.Range('E15').Value = '  -3.62%  '
$ws.Range('D16').Value = '0.0₅7941'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = '''64.52'
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').Value = '25.935.67'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '''1.002'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '''4.617'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').Value = '''184.57'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').Value = '''10.02'
$ws.Range('E22').Value = '  -4.09%  '
$ws.Range('D23').Value = '''6.069'
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').Value = '''1.003'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '''145.73'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('D26').Value = '''0.1207'
$ws.Range('E26').Value = '  -3.34%  '
$ws.Range('D27').Value = '''7.345'
$ws.Range('E27').Value = '  -2.99%  '
$ws.Range('D28').Value = '''15.50'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').Value = '''0.05935'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('E31').Value = '  -3.59%  '
$ws.Range('D32').Value = '''3.349'
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('D33').Value = '''3.354'
$ws.Range('E33').Value = '  -4.13%  '
$ws.Range('D34').Value = '''1.607'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').Value = '''0.9705'
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('D36').Value = '''2.382'
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('D37').Value = '''2.737'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = '''0.5778'
$ws.Range('E38').Value = '  -4.82%  '
$ws.Range('D39').Value = '''0.01596'
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('D41').Value = '''0.8386'
$ws.Range('E41').Value = '  -3.88%  '
$ws.Range('D42').Value = '''5.631'
$ws.Range('E42').Value = '  -7.65%  '
$ws.Range('D43').Value = '1.015.73'
$ws.Range('E43').Value = '  -6.29%  '
$ws.Range('D44').Value = '''99.55'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('D45').Value = '1.783.26'
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').Value = '0.0₈106'
$ws.Range('E46').Value = '  -3.51%  '
$ws.Range('D47').Value = '''0.9994'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '''54.25'
$ws.Range('E48').Value = '  -3.33%  '
$ws.Range('D49').Value = '''7.948'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').Value = '''0.05172'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('D51').Value = '''0.4225'
$ws.Range('E51').Value = '  -0.77%  '

# The leading apostrophe above forces Excel to keep these
# numeric-looking strings as literal text (matching the original
# inline-string cell type) instead of silently parsing them into
# numbers. Resetting the style back to Normal afterwards clears
# the 'quote prefix' flag that the apostrophe trick leaves behind,
# so the cell formatting is left exactly as it was originally.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
